$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.377249002456665
$ws.Range("B1").Value = 2.675835371017456
$ws.Range("C1").Value = 2.316540002822876
$ws.Range("D1").Value = 2.475037813186646
$ws.Range("E1").Value = 2.979245901107788
